# Penalty/Reward system update (unfinished) — shifts the forecast weeks
# forward by one and refreshes the forecast + summary numbers.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet --------------------------------------------
# Column B holds the week-start date as literal text (not a real date), so
# force a text number format before writing to keep it from being
# auto-converted into a date serial number.
$ws1.Range("B2:B17").NumberFormat = "@"

$ws1.Range("B2").Value  = "2025-01-12"
$ws1.Range("D2").Value  = 49

$ws1.Range("B3").Value  = "2025-01-19"
$ws1.Range("D3").Value  = 49

$ws1.Range("B4").Value  = "2025-01-26"
$ws1.Range("D4").Value  = 50

$ws1.Range("B5").Value  = "2025-02-02"
$ws1.Range("D5").Value  = 50

$ws1.Range("B6").Value  = "2025-02-09"
$ws1.Range("D6").Value  = 49

$ws1.Range("B7").Value  = "2025-02-16"
$ws1.Range("D7").Value  = 46

$ws1.Range("B8").Value  = "2025-02-23"
$ws1.Range("D8").Value  = 43

$ws1.Range("B9").Value  = "2025-03-02"
$ws1.Range("D9").Value  = 42

$ws1.Range("B10").Value = "2025-03-09"
$ws1.Range("D10").Value = 42

$ws1.Range("B11").Value = "2025-03-16"
$ws1.Range("D11").Value = 42

$ws1.Range("B12").Value = "2025-03-23"
$ws1.Range("D12").Value = 40

$ws1.Range("B13").Value = "2025-03-30"
$ws1.Range("D13").Value = 37

$ws1.Range("B14").Value = "2025-04-06"
$ws1.Range("D14").Value = 35

$ws1.Range("B15").Value = "2025-04-13"
$ws1.Range("D15").Value = 34

$ws1.Range("B16").Value = "2025-04-20"
$ws1.Range("D16").Value = 35

$ws1.Range("B17").Value = "2025-04-27"
$ws1.Range("D17").Value = 36

# --- Summary sheet ----------------------------------------------------------
$ws2.Range("B2:B15").NumberFormat = "@"

$ws2.Range("B2").Value  = "2023-01-29 to 2025-01-05"
$ws2.Range("B5").Value  = "42"
$ws2.Range("B6").Value  = "43"
$ws2.Range("B7").Value  = "21"
$ws2.Range("B8").Value  = "4036 units"
$ws2.Range("B9").Value  = "679"
$ws2.Range("B10").Value = "378"
$ws2.Range("B11").Value = "198"
$ws2.Range("B12").Value = "50"
$ws2.Range("B13").Value = "2025-01-26"
$ws2.Range("B14").Value = "34"
$ws2.Range("B15").Value = "2025-04-13"
